$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new key/value rows for test case / build tracking
# (values entered in this order to match shared-string allocation order)
$ws.Range("B10").Value = "PP_2_6_3_B2_P2"
$ws.Range("B9").Value = "Playwright_test_set"
$ws.Range("A9").Value = "feature_name"
$ws.Range("A10").Value = "build_version"

# Update the active cell selection to match the saved state
$ws.Range("F11").Select()
